# Update Ltbp1-Itgb5 LR-pair table to include "ECs" as a third sending
# cluster (in addition to the existing "FAPs" and "sCs"), per Dr Hou's
# advice. This expands the data from a 2x3 (sender x target) grid to a
# full 3x3 grid and refreshes all of the associated NATMI metrics.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> ECs (Ltbp1/Itgb5)
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Ltbp1"
$ws.Cells.Item(2, 3).Value = "Itgb5"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 2
$ws.Cells.Item(2, 6).Value = 0.6666666666666666
$ws.Cells.Item(2, 7).Value = 4.322843666666667
$ws.Cells.Item(2, 8).Value = 12.968531
$ws.Cells.Item(2, 9).Value = 0.06371657946635806
$ws.Cells.Item(2, 10).Value = 0.06371657946635806
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 8.309350333333333
$ws.Cells.Item(2, 14).Value = 24.928051
$ws.Cells.Item(2, 15).Value = 0.1535033474258946
$ws.Cells.Item(2, 16).Value = 0.1535033474258946
$ws.Cells.Item(2, 17).Value = 35.92002246256455
$ws.Cells.Item(2, 18).Value = 323.280202163081
$ws.Cells.Item(2, 19).Value = 0.009780708234613982
$ws.Cells.Item(2, 20).Value = 0.009780708234613984

# Row 3: ECs -> FAPs (Ltbp1/Itgb5)
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Ltbp1"
$ws.Cells.Item(3, 3).Value = "Itgb5"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 2
$ws.Cells.Item(3, 6).Value = 0.6666666666666666
$ws.Cells.Item(3, 7).Value = 4.322843666666667
$ws.Cells.Item(3, 8).Value = 12.968531
$ws.Cells.Item(3, 9).Value = 0.06371657946635806
$ws.Cells.Item(3, 10).Value = 0.06371657946635806
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 37.153391
$ws.Cells.Item(3, 14).Value = 111.460173
$ws.Cells.Item(3, 15).Value = 0.6863556906301786
$ws.Cells.Item(3, 16).Value = 0.6863556906301786
$ws.Cells.Item(3, 17).Value = 160.6083009795403
$ws.Cells.Item(3, 18).Value = 1445.474708815863
$ws.Cells.Item(3, 19).Value = 0.04373223690422485
$ws.Cells.Item(3, 20).Value = 0.04373223690422485

# Row 4: ECs -> sCs (Ltbp1/Itgb5)
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Ltbp1"
$ws.Cells.Item(4, 3).Value = "Itgb5"
$ws.Cells.Item(4, 4).Value = "sCs"
$ws.Cells.Item(4, 5).Value = 2
$ws.Cells.Item(4, 6).Value = 0.6666666666666666
$ws.Cells.Item(4, 7).Value = 4.322843666666667
$ws.Cells.Item(4, 8).Value = 12.968531
$ws.Cells.Item(4, 9).Value = 0.06371657946635806
$ws.Cells.Item(4, 10).Value = 0.06371657946635806
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 8.668653666666666
$ws.Cells.Item(4, 14).Value = 26.005961
$ws.Cells.Item(4, 15).Value = 0.1601409619439267
$ws.Cells.Item(4, 16).Value = 0.1601409619439267
$ws.Cells.Item(4, 17).Value = 37.47323460147678
$ws.Cells.Item(4, 18).Value = 337.259111413291
$ws.Cells.Item(4, 19).Value = 0.01020363432751923
$ws.Cells.Item(4, 20).Value = 0.01020363432751923

# Row 5: FAPs -> ECs (Ltbp1/Itgb5)
$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 2).Value = "Ltbp1"
$ws.Cells.Item(5, 3).Value = "Itgb5"
$ws.Cells.Item(5, 4).Value = "ECs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 29.080246
$ws.Cells.Item(5, 8).Value = 87.24073800000001
$ws.Cells.Item(5, 9).Value = 0.4286284557195201
$ws.Cells.Item(5, 10).Value = 0.4286284557195201
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 8.309350333333333
$ws.Cells.Item(5, 14).Value = 24.928051
$ws.Cells.Item(5, 15).Value = 0.1535033474258946
$ws.Cells.Item(5, 16).Value = 0.1535033474258946
$ws.Cells.Item(5, 17).Value = 241.6379517935153
$ws.Cells.Item(5, 18).Value = 2174.741566141638
$ws.Cells.Item(5, 19).Value = 0.06579590275493817
$ws.Cells.Item(5, 20).Value = 0.06579590275493817

# Row 6: FAPs -> FAPs (Ltbp1/Itgb5)
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Ltbp1"
$ws.Cells.Item(6, 3).Value = "Itgb5"
$ws.Cells.Item(6, 4).Value = "FAPs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 29.080246
$ws.Cells.Item(6, 8).Value = 87.24073800000001
$ws.Cells.Item(6, 9).Value = 0.4286284557195201
$ws.Cells.Item(6, 10).Value = 0.4286284557195201
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 37.153391
$ws.Cells.Item(6, 14).Value = 111.460173
$ws.Cells.Item(6, 15).Value = 0.6863556906301786
$ws.Cells.Item(6, 16).Value = 0.6863556906301786
$ws.Cells.Item(6, 17).Value = 1080.429750014186
$ws.Cells.Item(6, 18).Value = 9723.867750127674
$ws.Cells.Item(6, 19).Value = 0.2941915797491182
$ws.Cells.Item(6, 20).Value = 0.2941915797491181

# Row 7: FAPs -> sCs (Ltbp1/Itgb5)
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Ltbp1"
$ws.Cells.Item(7, 3).Value = "Itgb5"
$ws.Cells.Item(7, 4).Value = "sCs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 29.080246
$ws.Cells.Item(7, 8).Value = 87.24073800000001
$ws.Cells.Item(7, 9).Value = 0.4286284557195201
$ws.Cells.Item(7, 10).Value = 0.4286284557195201
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 8.668653666666666
$ws.Cells.Item(7, 14).Value = 26.005961
$ws.Cells.Item(7, 15).Value = 0.1601409619439267
$ws.Cells.Item(7, 16).Value = 0.1601409619439267
$ws.Cells.Item(7, 17).Value = 252.0865811154687
$ws.Cells.Item(7, 18).Value = 2268.779230039218
$ws.Cells.Item(7, 19).Value = 0.06864097321546377
$ws.Cells.Item(7, 20).Value = 0.06864097321546375

# Row 8: sCs -> ECs (Ltbp1/Itgb5)
$ws.Cells.Item(8, 1).Value = "sCs"
$ws.Cells.Item(8, 2).Value = "Ltbp1"
$ws.Cells.Item(8, 3).Value = "Itgb5"
$ws.Cells.Item(8, 4).Value = "ECs"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 34.44179
$ws.Cells.Item(8, 8).Value = 103.32537
$ws.Cells.Item(8, 9).Value = 0.5076549648141219
$ws.Cells.Item(8, 10).Value = 0.5076549648141219
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 8.309350333333333
$ws.Cells.Item(8, 14).Value = 24.928051
$ws.Cells.Item(8, 15).Value = 0.1535033474258946
$ws.Cells.Item(8, 16).Value = 0.1535033474258946
$ws.Cells.Item(8, 17).Value = 286.1888992170966
$ws.Cells.Item(8, 18).Value = 2575.70009295387
$ws.Cells.Item(8, 19).Value = 0.07792673643634244
$ws.Cells.Item(8, 20).Value = 0.07792673643634246

# Row 9: sCs -> FAPs (Ltbp1/Itgb5)
$ws.Cells.Item(9, 1).Value = "sCs"
$ws.Cells.Item(9, 2).Value = "Ltbp1"
$ws.Cells.Item(9, 3).Value = "Itgb5"
$ws.Cells.Item(9, 4).Value = "FAPs"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 34.44179
$ws.Cells.Item(9, 8).Value = 103.32537
$ws.Cells.Item(9, 9).Value = 0.5076549648141219
$ws.Cells.Item(9, 10).Value = 0.5076549648141219
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 37.153391
$ws.Cells.Item(9, 14).Value = 111.460173
$ws.Cells.Item(9, 15).Value = 0.6863556906301786
$ws.Cells.Item(9, 16).Value = 0.6863556906301786
$ws.Cells.Item(9, 17).Value = 1279.62929060989
$ws.Cells.Item(9, 18).Value = 11516.66361548901
$ws.Cells.Item(9, 19).Value = 0.3484318739768357
$ws.Cells.Item(9, 20).Value = 0.3484318739768357

# Row 10: sCs -> sCs (Ltbp1/Itgb5)
$ws.Cells.Item(10, 1).Value = "sCs"
$ws.Cells.Item(10, 2).Value = "Ltbp1"
$ws.Cells.Item(10, 3).Value = "Itgb5"
$ws.Cells.Item(10, 4).Value = "sCs"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 34.44179
$ws.Cells.Item(10, 8).Value = 103.32537
$ws.Cells.Item(10, 9).Value = 0.5076549648141219
$ws.Cells.Item(10, 10).Value = 0.5076549648141219
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 8.668653666666666
$ws.Cells.Item(10, 14).Value = 26.005961
$ws.Cells.Item(10, 15).Value = 0.1601409619439267
$ws.Cells.Item(10, 16).Value = 0.1601409619439267
$ws.Cells.Item(10, 17).Value = 298.5639491700633
$ws.Cells.Item(10, 18).Value = 2687.07554253057
$ws.Cells.Item(10, 19).Value = 0.08129635440094375
$ws.Cells.Item(10, 20).Value = 0.08129635440094375

